$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/source-record-id"

# Version: 7.0.0 -> 8.0.0
$wsMeta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet updates ---
$wsElements = $wb.Worksheets.Item("Elements")

# The "Fixed Value" of Extension.url (row 5, column Q) mirrors the same URL
# string as the Metadata sheet's URL property, so it must be kept in sync.
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/source-record-id"

# The root "Extension" row's Constraint(s) cell (AI2) no longer carries the
# ele-1/ext-1 invariant text in the regenerated IG output - it now only
# appears on the Extension.extension row (AI4, already populated).
$wsElements.Range("AI2").Value = ""
